$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-6.51%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.42%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.08%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07903"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.80%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.928"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-10.59%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.745"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.59%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.013"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.92%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.868"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.43%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9223"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.62%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1152"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "14.53%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1833"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.84%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09320"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.57%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03537"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.65%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09881"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.63%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001387"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.52%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005841"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.45%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.509"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.80%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3444"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.15%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1309"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.60%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.048"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.08%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2398"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.89%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04501"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.18%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001216"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.19%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004573"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.59%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001249"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.81%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.86%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01891"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.78%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04693"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.99%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007568"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.80%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009559"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "22.36%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1322"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.58%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002109"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.03%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01115"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.76%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006021"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-31.38%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.01%"
